# Auto-generated COM-interop script applying the cryptos.xlsx price/volume refresh
# described by the commit "Updated cryptos list ... with GitHub Actions".
# Updates Price (D) and Volume(1h) (E) text columns; four coin pairs also swapped
# ranking position, so their Coin (B) / Link (C) cells are rewritten too.
#
# All Price values are plain text in the sheet (e.g. "236.16", "3.572.68" using
# dot-grouped thousands) so we prefix a leading apostrophe when assigning through
# .Value to stop Excel from auto-coercing the text into a Number/date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'95.107.96"
$ws.Range("E2").Value = "  -1.20%  "

$ws.Range("D3").Value = "'3.572.68"

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'236.16"
$ws.Range("E5").Value = "  -1.62%  "

$ws.Range("D6").Value = "'656.69"
$ws.Range("E6").Value = "  +2.35%  "

$ws.Range("E7").Value = "  -0.63%  "

$ws.Range("E8").Value = "  -0.27%  "

$ws.Range("E9").Value = "  +0.10%  "

$ws.Range("E10").Value = "  +0.01%  "

$ws.Range("D11").Value = "'3.570.75"
$ws.Range("E11").Value = "  -1.52%  "

$ws.Range("E12").Value = "  +0.96%  "

$ws.Range("D13").Value = "'42.38"
$ws.Range("E13").Value = "  -2.43%  "

$ws.Range("D14").Value = "'6.46"
$ws.Range("E14").Value = "  +1.90%  "

$ws.Range("D15").Value = "'4.236.10"
$ws.Range("E15").Value = "  -1.88%  "

$ws.Range("D16").Value = "'95.009.67"
$ws.Range("E16").Value = "  -1.20%  "

$ws.Range("E17").Value = "  -0.35%  "

$ws.Range("D18").Value = "'3.569.95"
$ws.Range("E18").Value = "  -1.45%  "

$ws.Range("D19").Value = "'7.75"
$ws.Range("E19").Value = "  -7.46%  "

$ws.Range("D20").Value = "'12.65"
$ws.Range("E20").Value = "  -5.13%  "

$ws.Range("D21").Value = "'17.87"
$ws.Range("E21").Value = "  -2.36%  "

$ws.Range("E22").Value = "  +0.73%  "

$ws.Range("B23").Value = "Stellar"
$ws.Range("C23").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D23").Value = "'0.483"
$ws.Range("E23").Value = "  -2.92%  "

$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").Value = "'509.03"
$ws.Range("E24").Value = "  -1.60%  "

$ws.Range("E25").Value = "  +4.15%  "

$ws.Range("E26").Value = "  -0.09%  "

$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D27").Value = "'95.04"
$ws.Range("E27").Value = "  -2.31%  "

$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").Value = "'12.88"
$ws.Range("E28").Value = "  +3.00%  "

$ws.Range("D29").Value = "'3.764.21"
$ws.Range("E29").Value = "  -1.59%  "

$ws.Range("E30").Value = "  -2.77%  "

$ws.Range("D31").Value = "'0.145"
$ws.Range("E31").Value = "  +1.36%  "

$ws.Range("D32").Value = "'11.56"
$ws.Range("E32").Value = "  -0.02%  "

$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("D34").Value = "'0.997"
$ws.Range("E34").Value = "  +0.81%  "

$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").Value = "'32.08"
$ws.Range("E35").Value = "  +2.74%  "

$ws.Range("B36").Value = "Cronos"
$ws.Range("C36").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D36").Value = "'0.177"
$ws.Range("E36").Value = "  -2.51%  "

$ws.Range("D37").Value = "'1.72"
$ws.Range("E37").Value = "  +16.69%  "

$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").Value = "'0.561"
$ws.Range("E38").Value = "  -2.18%  "

$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "'603.39"
$ws.Range("E39").Value = "  +5.85%  "

$ws.Range("D40").Value = "'8.58"
$ws.Range("E40").Value = "  +9.30%  "

$ws.Range("E41").Value = "  +0.14%  "

$ws.Range("E42").Value = "  +0.12%  "

$ws.Range("E43").Value = "  -3.25%  "

$ws.Range("D44").Value = "'1.86"
$ws.Range("E44").Value = "  +7.17%  "

$ws.Range("D45").Value = "'35.11"
$ws.Range("E45").Value = "  +24.53%  "

$ws.Range("D46").Value = "'5.74"
$ws.Range("E46").Value = "  +0.08%  "

$ws.Range("D47").Value = "'2.30"
$ws.Range("E47").Value = "  +3.97%  "

$ws.Range("E49").Value = "  -3.26%  "

$ws.Range("D50").Value = "'3.50"
$ws.Range("E50").Value = "  +0.03%  "

$ws.Range("D51").Value = "'8.21"
$ws.Range("E51").Value = "  +0.40%  "
